$wb = $excel.ActiveWorkbook

# ALC!row5
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 118.42857
$ws.Range("I5").Value = 100.61539
$ws.Range("K5").Value = 100.61539
$ws.Range("M5").Value = 14.38461

# ALC!row9
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 1492.2222
$ws.Range("I9").Value = 1492.2222
$ws.Range("K9").Value = 1492.2222
$ws.Range("M9").Value = -1323.2222

# ALC!row28
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 42.25
$ws.Range("I28").Value = 42.25
$ws.Range("K28").Value = 42.25
$ws.Range("M28").Value = 442.75

# ALC!row33
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("M33").Value = ""
$ws.Range("H33").Value = 400
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0

# ALC!row76
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 4996.5
$ws.Range("I76").Value = 4993
$ws.Range("J76").Value = 5000
$ws.Range("K76").Value = 4993
$ws.Range("L76").Value = 5000
$ws.Range("M76").Value = -4678
$ws.Range("N76").Value = -5630

# ALC!row79
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 4996.5
$ws.Range("I79").Value = 4993
$ws.Range("J79").Value = 5000
$ws.Range("K79").Value = 4993
$ws.Range("L79").Value = 5000
$ws.Range("M79").Value = -3901
$ws.Range("N79").Value = -7184

# ALC!row86
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 1779.8
$ws.Range("I86").Value = 1733
$ws.Range("J86").Value = 1850
$ws.Range("K86").Value = 1733
$ws.Range("L86").Value = 1850
$ws.Range("M86").Value = -610
$ws.Range("N86").Value = -4096

# ALC!row88
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 5080
$ws.Range("I88").Value = 4800
$ws.Range("J88").Value = 5500
$ws.Range("K88").Value = 4800
$ws.Range("L88").Value = 5500
$ws.Range("M88").Value = -4394
$ws.Range("N88").Value = -6312

# ALC!row89
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 1779.8
$ws.Range("I89").Value = 1733
$ws.Range("J89").Value = 1850
$ws.Range("K89").Value = 8665
$ws.Range("L89").Value = 9250
$ws.Range("M89").Value = -3049
$ws.Range("N89").Value = -20482

# ALC!row91
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H91").Value = 5080
$ws.Range("I91").Value = 4800
$ws.Range("J91").Value = 5500
$ws.Range("K91").Value = 4800
$ws.Range("L91").Value = 5500
$ws.Range("M91").Value = -3396
$ws.Range("N91").Value = -8308

# ALC!row138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 6270.4287
$ws.Range("I138").Value = 4973.25
$ws.Range("K138").Value = 14919.75
$ws.Range("M138").Value = -9779.75

# ARM!row63
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3449.75
$ws.Range("I63").Value = 3449.75
$ws.Range("K63").Value = 3449.75
$ws.Range("M63").Value = -2763.75

# ARM!row66
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 3449.75
$ws.Range("I66").Value = 3449.75
$ws.Range("K66").Value = 17248.75
$ws.Range("M66").Value = -13816.75

# ARM!row80
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H80").Value = 1689666.6
$ws.Range("I80").Value = 34500
$ws.Range("K80").Value = 34500
$ws.Range("M80").Value = -33502

# ARM!row83
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H83").Value = 1689666.6
$ws.Range("I83").Value = 34500
$ws.Range("K83").Value = 103500
$ws.Range("M83").Value = -98508

# ARM!row97
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 613.25
$ws.Range("I97").Value = 504.42856
$ws.Range("K97").Value = 504.42856
$ws.Range("M97").Value = -8.428560000000004

# BSM!row86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4749.826
$ws.Range("I86").Value = 2741.5
$ws.Range("J86").Value = 5820.933
$ws.Range("K86").Value = 2741.5
$ws.Range("L86").Value = 5820.933
$ws.Range("M86").Value = -1618.5
$ws.Range("N86").Value = -8066.933

# BSM!row89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 4749.826
$ws.Range("I89").Value = 2741.5
$ws.Range("J89").Value = 5820.933
$ws.Range("K89").Value = 13707.5
$ws.Range("L89").Value = 29104.665
$ws.Range("M89").Value = -8091.5
$ws.Range("N89").Value = -40336.665

# BSM!row94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("N94").Value = ""
$ws.Range("H94").Value = 1000
$ws.Range("I94").Value = 1000
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 1000
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -549

# BSM!row99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1199.6666
$ws.Range("I99").Value = 1199.6666
$ws.Range("K99").Value = 1199.6666
$ws.Range("M99").Value = 298.3334

# BSM!row107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 6596.5557
$ws.Range("I107").Value = 6983.625
$ws.Range("K107").Value = 6983.625
$ws.Range("M107").Value = -5063.625

# CRP!row31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 37848.625
$ws.Range("I31").Value = 2198.25
$ws.Range("K31").Value = 2198.25
$ws.Range("M31").Value = -1903.25

# CRP!row34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 37848.625
$ws.Range("I34").Value = 2198.25
$ws.Range("K34").Value = 2198.25
$ws.Range("M34").Value = -1996.25

# CRP!row86
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 7213.4287
$ws.Range("J86").Value = 999
$ws.Range("L86").Value = 999
$ws.Range("N86").Value = -3245

# CRP!row89
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 7213.4287
$ws.Range("J89").Value = 999
$ws.Range("L89").Value = 4995
$ws.Range("N89").Value = -16227

# CRP!row105
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 2090.2307
$ws.Range("I105").Value = 2106.6365
$ws.Range("J105").Value = 2000
$ws.Range("K105").Value = 2106.6365
$ws.Range("L105").Value = 2000
$ws.Range("M105").Value = -359.6365000000001
$ws.Range("N105").Value = -5494

# CRP!row107
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 676.5333000000001
$ws.Range("I107").Value = 656.5714
$ws.Range("J107").Value = 694
$ws.Range("K107").Value = 656.5714
$ws.Range("L107").Value = 694
$ws.Range("M107").Value = 1263.4286
$ws.Range("N107").Value = -4534

# CRP!row132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3767.5
$ws.Range("I132").Value = 3630
$ws.Range("K132").Value = 10890
$ws.Range("M132").Value = -8360

# CUL!row34
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1586.75
$ws.Range("J34").Value = 2115
$ws.Range("L34").Value = 6345
$ws.Range("N34").Value = -6513

# CUL!row39
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 745.8333
$ws.Range("J39").Value = 725
$ws.Range("L39").Value = 2175
$ws.Range("N39").Value = -2763

# CUL!row55
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 7000
$ws.Range("J55").Value = 7000
$ws.Range("L55").Value = 21000
$ws.Range("N55").Value = -21354

# CUL!row112
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H112").Value = 3016.6667
$ws.Range("I112").Value = 3016.6667
$ws.Range("K112").Value = 9050.000100000001
$ws.Range("M112").Value = -7942.000100000001

# CUL!row130
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H130").Value = 1500
$ws.Range("I130").Value = 1500
$ws.Range("K130").Value = 4500
$ws.Range("M130").Value = 520

# CUL!row137
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("N137").Value = ""
$ws.Range("H137").Value = 4930
$ws.Range("I137").Value = 4930
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 14790
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -9690

# GSM!row97
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2015.6
$ws.Range("I97").Value = 1769.7
$ws.Range("K97").Value = 1769.7
$ws.Range("M97").Value = -1273.7

# LTW!row82
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2366.3333
$ws.Range("I82").Value = 2474.625
$ws.Range("K82").Value = 2474.625
$ws.Range("M82").Value = -2113.625

# LTW!row85
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 2366.3333
$ws.Range("I85").Value = 2474.625
$ws.Range("K85").Value = 2474.625
$ws.Range("M85").Value = -1226.625

# LTW!row133
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H133").Value = 71775.336
$ws.Range("J133").Value = 71775.336
$ws.Range("L133").Value = 71775.336
$ws.Range("N133").Value = -76835.336

# WVR!row30
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("N30").Value = ""
$ws.Range("H30").Value = 15000
$ws.Range("I30").Value = 15000
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 15000
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = -14893

# WVR!row100
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 582.1111
$ws.Range("I100").Value = 404.5
$ws.Range("K100").Value = 809
$ws.Range("M100").Value = -268

# WVR!row109
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 42671
$ws.Range("I109").Value = 40342
$ws.Range("J109").Value = 45000
$ws.Range("K109").Value = 40342
$ws.Range("L109").Value = 45000
$ws.Range("M109").Value = -38955
$ws.Range("N109").Value = -47774
